$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowCH($row, $vals) {
    $arr = New-Object "object[,]" 1,6
    for ($i = 0; $i -lt 6; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("C$row`:H$row").Value2 = $arr
}

Set-RowCH 2 @(-3.373677730560303, 8.367032051086426, -1.811180233955384, 0.08162501163598926, -0.581960884536185, 0.1855354215495473)
Set-RowCH 3 @(-3.436809062957764, 8.002476692199707, -2.409347295761108, 0.6193056222869129, 0.1286169377768912, -0.6118523001489119)
Set-RowCH 4 @(-4.233324527740479, 7.347033500671387, -3.329159736633301, 0.8970928599194783, 0.8833782091373344, -3.204748107165821)
Set-RowCH 5 @(-4.557640552520752, 5.48013973236084, -4.308228492736816, -0.161268437780985, 0.124489894727381, -2.472156795059762)
Set-RowCH 6 @(-6.2513108253479, 4.66443920135498, -5.362371444702148, -2.366607433412133, 2.737767719640964, -1.411724490363423)
Set-RowCH 7 @(-6.407595157623291, 5.533173084259033, -4.768703460693359, -2.001652684517023, 3.962463135399469, -1.871139858554049)
Set-RowCH 8 @(-6.354794025421143, 3.896303415298462, -3.691534996032715, 0.06661038346043445, 0.1919569562121135, 0.05933587689225263)
Set-RowCH 9 @(55.61550140380859, -7.603207588195801, 15.50844955444336, -0.02138400813791778, 0.1378432366906143, -0.02639385003869168)
Set-RowCH 10 @(-3.817588806152344, 5.203746795654297, -4.701224803924561, 0.1729754304013602, -0.1063352536137511, -0.01823656188278665)
Set-RowCH 11 @(2.208892345428467, 4.523828029632568, -5.306243419647217, 0.4490567220420372, -0.3232520183593762, -0.09577174720967686)
Set-RowCH 12 @(4.054627418518066, 6.320723056793213, -5.26066255569458, 0.04528970053283177, 0.1171147901473975, 0.09616285126383706)
Set-RowCH 13 @(3.682106018066406, 6.276570796966553, -4.233500480651856, -0.0392033953855677, 0.2259382031312803, 0.1824885385065544)
Set-RowCH 14 @(4.067971229553223, 7.197247505187988, -6.232089042663574, -0.03262541961015729, -0.1276075131282566, 0.02622995907213681)
Set-RowCH 15 @(3.399141788482666, 6.60714864730835, -4.322453498840332, 0.003914752202790014, -0.2204664746617399, -0.05295903471911827)
Set-RowCH 16 @(2.337307453155518, 7.003424167633057, -5.783712863922119, 0.01699993268745699, 0.00523332715397923, 0.05201666502327448)
Set-RowCH 17 @(2.072117805480957, 7.399670124053955, -5.612285137176514, -0.01029531448716063, -0.03043524366689888, -0.04977434319330424)
Set-RowCH 18 @(2.665160179138184, 6.522679805755615, -4.401774883270264, -0.002801040596351329, -0.03519552141973155, -0.04244768401471578)
Set-RowCH 19 @(3.009289264678955, 6.934535980224609, -4.890542030334473, 0.006436434000866758, 0.02465437250438991, 0.02476239040857406)
Set-RowCH 20 @(3.372682571411133, 7.015328407287598, -4.779996871948242, 0.0004395249685863074, 0.0006667370567233805, 0.01120416263527254)
Set-RowCH 21 @(3.250783443450928, 6.881857395172119, -4.688606262207031, -0.004521892917137975, -0.01801680073868935, -0.009371567176791164)

# Row 22 (old last data row) is no longer part of the dataset.
$ws.Rows.Item(22).Delete()
